# Generate Report for Handback
# Update timestamps / status recorded on the handback status report.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" column (G) for the
# 16f05ebc-... row (row 2) and the 62c8e010-... row (row 4) shared the
# same generated timestamp; bump it forward to the newly generated value.
$wsOverview.Range("G2").Value = "2016-09-06 02:18:21"
$wsOverview.Range("G4").Value = "2016-09-06 02:18:21"

# zh-cn sheet: Priority changed from "ht" to "mt" for both rows that used it,
# and the handoff/handback datetimes for the 16f05ebc-... file were refreshed.
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("E4").Value = "mt"

$wsZhCn.Range("H2").Value = "2016-09-06 02:18:16"
$wsZhCn.Range("H4").Value = "2016-09-06 02:18:16"

$wsZhCn.Range("K2").Value = "2016-09-06 02:18:34"
$wsZhCn.Range("K4").Value = "2016-09-06 02:18:34"

# de-de sheet: same priority update and a refreshed handback datetime for
# the 16f05ebc-... file.
$wsDeDe.Range("E2").Value = "mt"
$wsDeDe.Range("E4").Value = "mt"

$wsDeDe.Range("K2").Value = "2016-09-06 02:18:41"
$wsDeDe.Range("K4").Value = "2016-09-06 02:18:41"
